$d = $word.ActiveDocument

# 1. Replace the lead-in text before the hyperlink:
#      " As Geotechnologist in Residence for "
#    with the new bio lead-in:
#      " He is co-founder and director of Malstow Geospatial, a consultancy
#        firm offering bespoke consulting and services in the geospatial,
#        geotechnology, maps and location based services fields. This means
#        Gary is currently consulting as Head of APIs for the "
$lead = $d.Content
$leadFound = $lead.Find.Execute(" As Geotechnologist in Residence for ")
if ($leadFound) {
    $lead.Text = " He is co-founder and director of Malstow Geospatial, a consultancy firm offering bespoke consulting and services in the geospatial, geotechnology, maps and location based services fields. This means Gary is currently consulting as Head of APIs for the "
}

# 2. Update the hyperlink's visible text from "Lokku" to "Ordnance Survey"
#    (use TextToDisplay so the hyperlink style/target, i.e. r:id/rId5, are
#    preserved unchanged, matching the diff which leaves the relationship
#    target untouched).
for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    if ($h.TextToDisplay -eq "Lokku") {
        $h.TextToDisplay = "Ordnance Survey"
    }
}

# 3. Replace the trailing sentence after the hyperlink:
#      ", Gary is helping to advance open geospatial technologies and bring
#        them to new markets."
#    with:
#      ", the United Kingdom's national mapping agency."
# The old text starts right where the hyperlink ends, so a plain
# Find/Replace (or a direct Range.Text= on the found range) would make the
# new run inherit the hyperlink's run style. Instead, anchor the insertion
# at the END of the old range -- which borrows the plain formatting of the
# " A Fellow of the " text that follows -- and then delete the stale text
# that is left behind in front of it.
$old = $d.Content
$oldFound = $old.Find.Execute(", Gary is helping to advance open geospatial technologies and bring them to new markets.")
if ($oldFound) {
    $oldStart = $old.Start
    $oldEnd = $old.End

    $insertPoint = $d.Range($oldEnd, $oldEnd)
    $insertPoint.InsertBefore(", the United Kingdom's national mapping agency.")

    $stale = $d.Range($oldStart, $oldEnd)
    $stale.Delete()
}
